# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
# -------------------------------------------------------------
# This NIT-9008983067 account-statement sheet gets:
#   1. A new overdue-period row for worker "DORIA RUTH VEJARANO PARDO"
#      (same CC 51561138) for period 2509, inserted right after her
#      existing 2104 row (new row 18).
#   2. All the other workers' "Periodo Mora" updated from 2508 -> 2509.
#   3. The "VALOR MORA" total (E11) bumped by the new row's amount.
# -------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at 18 (pushes the rest of the table, and the
#    signature block further down, from 26/27 to 27/28).
$ws.Rows.Item(18).Insert()

# Pull the formatting (fonts/borders/number formats) from row 16 - the
# existing row for this same worker - so the new row matches the rest
# of the table exactly.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's data: same worker (CC 51561138, DORIA RUTH
# VEJARANO PARDO) but a different overdue period and amount.
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "51561138"
$ws.Range("D18").Value = "DORIA RUTH VEJARANO PARDO"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 560000
$ws.Range("G18").Value = 14000000

# 2) Update the remaining workers' Periodo Mora from 2508 to 2509
#    (rows shifted down by the insert above).
$ws.Range("E17").Value = "2509"
$ws.Range("E19").Value = "2509"
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"
$ws.Range("E22").Value = "2509"

# 3) Update the VALOR MORA total to reflect the new row
#    (621433 + 560000 = 1181433).
$ws.Range("E11").Value = 1181433

Write-Output "Edit complete"
